# Generate Report for Handoff
# Updates status from "In Translation" to "Ready for handoff" and refreshes the
# handoff timestamps on the Overview, zh-cn and de-de worksheets, then widens
# the Status columns so the new text fits.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Target column width (post-handoff "Status" columns) is 17.2159881591797
# characters in the saved OOXML. The COM column-width setter here only lands
# on a discrete (k+5)/6 pixel grid, so 16.333333333333332 is the closest
# input that rounds to the nearest reachable grid value (17.166666666666668).
$statusColWidth = 16.333333333333332

# --- Overview sheet ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-31 21:14:11"
$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-31 21:14:05"
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-31 21:14:11"
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
